$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript three char used in PEPE price (U+2083)
$sub3 = [char]0x2083

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.227.15'
$ws.Range("E2").Value = '  -0.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.426.63'
$ws.Range("E3").Value = '  -0.30%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.81'
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.75'
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.590'
$ws.Range("E8").Value = '  +1.93%  '
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.69'
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.356'
$ws.Range("E11").Value = '  -1.27%  '
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '25.14'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.859.31'
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.147.01'
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.411.96'
$ws.Range("E17").Value = '  +1.23%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.35'
$ws.Range("E18").Value = '  -0.76%  '
$ws.Range("E19").Value = '  +1.46%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '329.92'
$ws.Range("E20").Value = '  -1.51%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.75'
$ws.Range("E21").Value = '  -1.38%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.89'
$ws.Range("E23").Value = '  +0.67%  '
$ws.Range("E24").Value = '  +3.83%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.74'
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.15%  '
$ws.Range("E27").Value = '  +2.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = [string]::Concat('0.0', $sub3, '0778')
$ws.Range("E28").Value = '  -1.42%  '
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '170.53'
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.15'
$ws.Range("E31").Value = '  -2.81%  '
$ws.Range("E32").Value = '  +4.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.406'
$ws.Range("E33").Value = '  -3.33%  '
$ws.Range("E34").Value = '  -0.87%  '
$ws.Range("E35").Value = '  +1.74%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.25'
$ws.Range("E37").Value = '  +0.37%  '
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '334.21'
$ws.Range("E39").Value = '  +3.78%  '
$ws.Range("E40").Value = '  -0.77%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.91'
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '144.72'
$ws.Range("E42").Value = '  +1.76%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.68'
$ws.Range("E43").Value = '  -1.28%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.15'
$ws.Range("E44").Value = '  +2.55%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0969'
$ws.Range("E45").Value = '  +0.73%  '
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0518'
$ws.Range("E46").Value = '  -1.44%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.579'
$ws.Range("E47").Value = '  +0.89%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0224'
$ws.Range("E48").Value = '  -1.25%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.04'
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.59'
$ws.Range("E50").Value = '  -2.42%  '
$ws.Range("B51").Value = 'ZEEBU'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B5-YKN_zB+zeebu-zbu'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.67'
$ws.Range("E51").Value = '  -1.02%  '
